# xinyangapitesting own.xlsx — "all first three test cases run succesfully"
#
# The commit replaces the expected_result text used by the later login
# test cases (login_002..login_008) on the "login" sheet: the Chinese
# error message {"text":"账号或密码不正确"} is swapped for its English
# translation {"text":"Incorrect account or password"}. That string was
# only referenced from column G (expected_result), rows 3-9, of the
# "login" sheet, so writing the new literal value there reproduces the
# shared-string-table change (old string dropped, new one appended) as
# a natural side effect of saving.
#
# The rest of the diff (row/col reindexing inside sharedStrings.xml,
# the shifted <v> indices on every sheet, the dropped "宋体" font/xf
# entry) is a mechanical consequence of that single content edit and
# needs no separate action.
#
# A few purely cosmetic, non-data view/format tweaks from the diff are
# also reproduced where the object model exposes them: the selected
# cell on each sheet, and the widened "request_parameter"/"expected_result"
# columns (F/G) on the "login" sheet.

$wb = $excel.ActiveWorkbook

# --- content edit -----------------------------------------------------
$wsLogin = $wb.Worksheets.Item("login")
for ($r = 3; $r -le 9; $r++) {
    $wsLogin.Cells.Item($r, 7).Value = '{"text":"Incorrect account or password"}'
}

# --- cosmetic: column widths on the login sheet (F -> 43.5, G -> 31.5) -
# ColumnWidth's stored width runs ~5/6 wider than the value assigned
# (Excel's character-to-pixel rounding), so back the input off by 5/6
# to land exactly on the target stored widths.
$wsLogin.Columns.Item(6).ColumnWidth = 43.5 - (5 / 6)
$wsLogin.Columns.Item(7).ColumnWidth = 31.5 - (5 / 6)

# --- cosmetic: selection / active cell per sheet -----------------------
$wsRegister = $wb.Worksheets.Item("test_register")
$wsRegister.Range("E2").Select()

$wsFlow = $wb.Worksheets.Item("test_business_flow")
$wsFlow.Range("G11").Select()

$wsLogin.Activate()
$wsLogin.Range("G9").Select()
